$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64-84 down to rows 65-85
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with its data
$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44900
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = 300000000
$ws.Cells.Item(64, 7).Value = "Espárragos"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 1500
$ws.Cells.Item(64, 12).Value = 1500
$ws.Cells.Item(64, 13).Value = 1500
$ws.Cells.Item(64, 14).Value = "$/kilo"
$ws.Cells.Item(64, 15).Value = "Provincia de Linares"
$ws.Cells.Item(64, 16).Value = 1500
$ws.Cells.Item(64, 17).Value = 1
$ws.Cells.Item(64, 18).Value = "Hortaliza"
